$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A13").Value = "Atilius Fortunatianus"
$ws.Range("A14").Select()
